$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text-formatted numbers (dotted thousands
# separators, e.g. "29.334.64"); force text format before assigning so
# Excel does not reinterpret them as numeric/date values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.334.64'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.861.38'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7015'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.73'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07790'
$ws.Range("E8").Value = '  -2.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3044'
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("E10").Value = '  +6.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08156'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.838.40'
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.212'
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7145'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.13'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.317.33'
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.785'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.08'
$ws.Range("E18").Value = '  +2.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007764'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.15'
$ws.Range("E20").Value = '  -1.51%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.105.82'
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.509'
$ws.Range("E24").Value = '  +1.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.38'
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.884'
$ws.Range("E26").Value = '  -0.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1430'
$ws.Range("E27").Value = '  -1.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.06'
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.898'
$ws.Range("E29").Value = '  -4.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.371'
$ws.Range("E30").Value = '  -4.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.471'
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.291'
$ws.Range("E32").Value = '  -2.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.024'
$ws.Range("E33").Value = '  -0.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05153'
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.179'
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7042'
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9975'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.677'
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01839'
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.689'
$ws.Range("E40").Value = '  -1.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.176.42'
$ws.Range("E41").Value = '  +2.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9149'
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.008'
$ws.Range("E43").Value = '  +1.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.45'
$ws.Range("E44").Value = '  +1.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4232'
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.51'
$ws.Range("E47").Value = '  -1.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5352'
$ws.Range("E48").Value = '  -1.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.744'
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.123'
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.936'
$ws.Range("E51").Value = '  +0.28%  '
